$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the label-cell formatting (bold, bordered, centered) from A1 onto A2/A3
# before we rewrite the sheet contents.
$ws.Range("A1").Copy()
$ws.Range("A2:A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Columns C and D are being removed entirely from the table
$ws.Range("C1:D3").Clear()

# B column holds plain (unstyled) numbers - strip the header style B1 had
$ws.Range("B1:B3").Style = "Normal"

# Row 1
$ws.Range("A1").Value = "Puntuación promedio de los airbnb de la zona de Roberto y Clara"
$ws.Range("B1").Value = 3.3

# Row 2
$ws.Range("A2").Value = "Reseñas promedio de los airbnb de la zona de Roberto y Clara"
$ws.Range("B2").Value = 28

# Row 3
$ws.Range("A3").Value = "Precio promedio de los airbnb de la zona de Roberto y Clara"
$ws.Range("B3").Value = 83.59999999999999
